$d = $word.ActiveDocument

# Merge the title runs into a single run: "Client Meeting" + " 2" + " - AGMeeting"
# -> "Client Meeting 2 - AGMeeting"
$d.Content.Find.Execute("Client Meeting 2 " + [char]8211 + " AGMeeting", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Client Meeting 2 " + [char]8211 + " AGMeeting", 2)

# Update the meeting date text: "Oct.10/21" -> "Oct.1/21"
$d.Content.Find.Execute("Oct.10/21", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Oct.1/21", 2)

# Merge the "End Time" trailing runs: " " + "am pst" -> " am pst"
# (keep the "13" run, which carries its own rsid, untouched; only the
# space run and the "am pst" run - which share identical formatting -
# collapse into one run, matching the target XML.) Locate "am pst"
# dynamically (rather than hard-coded offsets), then replace just that
# found run in place so it merges with its immediately preceding,
# identically-formatted space run (but not the differently-rsid'd "13"
# run before it).
$endTimePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "End Time:*") {
        $endTimePara = $p
        break
    }
}

$locate = $endTimePara.Range.Duplicate
$locate.Find.Execute("am pst", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
$amRange = $d.Range($locate.Start, $locate.End)
$amRange.Find.Execute("am pst", $true, $false, $false, $false, $false,
                       $true, 1, $false, "am pst", 2)
